$wb = $excel.ActiveWorkbook

$sliders = $wb.Worksheets.Item("Sliders")
$commands = $wb.Worksheets.Item("Commands")

# --- Commands sheet: insert a new row for the "adjustSV" Artisan Command ---
# New row is inserted right above the existing "pidSV(<float>)" row (row 75),
# pushing it (and everything below) down by one row.
$commands.Rows.Item(75).Insert()
$commands.Rows.Item(75).RowHeight = 13.8

$commands.Range("B75").Value = "adjustSV(<int>)"
$commands.Range("C75").Value = "increases or decreases the current target SV value by <int>"

# The pre-existing pidSV row (now shifted to row 76) changes its parameter
# type from <float> to <int>.
$commands.Range("B76").Value = "pidSV(<int>)"

# --- Restore/update the view selection state on both sheets ---
$sliders.Activate()
$sliders.Range("B6").Select()

$commands.Activate()
$commands.Range("B76").Select()
